{"js": "// Bug Fixes and almost complete Error Example file\n//\n// 1) In the \"Compiling\" section, extend the compile command so it also\n//    lists the extra source files (Lexical_Analyzer.c and VM.c):\n//       \"To compile please run gcc main.c. \"\n//    becomes\n//       \"To compile please run gcc main.c Lexical_Analyzer.c VM.c. \"\n//\n// 2) The \"_GoBack\" bookmark (previously sitting at the very end of the\n//    \"Example: ./a.out ...\" paragraph) is relocated to sit right after the\n//    newly-typed text in the \"Compiling\" paragraph (this is exactly what\n//    Word does automatically: it always keeps \"_GoBack\" at the location of\n//    the most recent edit).\n\nconst body = context.document.body;\n\n// --- Step 1: insert the two extra filenames right after \"gcc main.c\" ---\nconst compileMatches = body.search(\"gcc main.c\", { matchCase: true, matchWholeWord: false });\ncompileMatches.load(\"items\");\nawait context.sync();\n\nif (compileMatches.items.length === 0) {\n  throw new Error(\"Could not locate 'gcc main.c' text to edit.\");\n}\n\nconst insertionRange = compileMatches.items[0].getRange(\"End\");\ninsertionRange.insertText(\" Lexical_Analyzer.c VM.c\", \"Start\");\nawait context.sync();\n\n// --- Step 2: move the \"_GoBack\" bookmark to just after the text we added ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst newFileMatches = body.search(\"VM.c\", { matchCase: true, matchWholeWord: false });\nnewFileMatches.load(\"items\");\nawait context.sync();\n\nif (newFileMatches.items.length === 0) {\n  throw new Error(\"Could not locate newly inserted 'VM.c' text for the bookmark.\");\n}\n\nconst bookmarkRange = newFileMatches.items[0].getRange(\"End\");\nbookmarkRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Bug Fixes and almost complete Error Example file\n#\n# 1) In the \"Compiling\" section, extend the compile command so it also\n#    lists the extra source files (Lexical_Analyzer.c and VM.c):\n#       \"To compile please run gcc main.c. \"\n#    becomes\n#       \"To compile please run gcc main.c Lexical_Analyzer.c VM.c. \"\n#\n# 2) The \"_GoBack\" bookmark (previously sitting at the very end of the\n#    \"Example: ./a.out ...\" paragraph) is relocated to sit right after the\n#    newly-typed text in the \"Compiling\" paragraph (this is exactly what\n#    Word does automatically: it always keeps \"_GoBack\" at the location of\n#    the most recent edit).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: insert the two extra filenames right after \"gcc main.c\" ---\n$editRange = $d.Content\n$found = $editRange.Find.Execute(\"gcc main.c\")\nif (-not $found) {\n    throw \"Could not locate 'gcc main.c' text to edit.\"\n}\n$editRange.Collapse(0)  # wdCollapseEnd\n$editRange.InsertAfter(\" Lexical_Analyzer.c VM.c\")\n\n# --- Step 2: move the \"_GoBack\" bookmark to just after the text we added ---\n$oldBookmark = $d.Bookmarks(\"_GoBack\")\n$oldBookmark.Delete()\n\n$bookmarkRange = $d.Content\n$foundVm = $bookmarkRange.Find.Execute(\"VM.c\")\nif (-not $foundVm) {\n    throw \"Could not locate newly inserted 'VM.c' text for the bookmark.\"\n}\n$bookmarkRange.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
